$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that flip COMPLETE (column D) from FALSE to TRUE, no notes added.
$completedRows = @(13, 16, 17, 25, 26, 27, 29, 32, 37)
foreach ($r in $completedRows) {
    $ws.Range("D$r").Value = $true
}

# Rows that flip COMPLETE to TRUE *and* get a reviewer note in column E.
$ws.Range("D22").Value = $true
$ws.Range("E22").Value = "(added into experiments section intro)"

$ws.Range("D23").Value = $true
$ws.Range("E23").Value = "(explained in the experiments section)"

$ws.Range("D24").Value = $true
$ws.Range("E24").Value = "It is in the methodology section, also I explained my rationale behind chosing computational complexity over simple O-analysis"
$ws.Rows.Item(24).RowHeight = 43.2

$ws.Range("D30").Value = $true
$ws.Range("E30").Value = "alread done in the ""Performance Comparison"" section"

# Recalculate so the shared F-column formulas and the D57/F57 summary update.
$excel.Calculate()

# Move the active selection back up to D3, which also drops the stale
# scrolled-down topLeftCell recorded in the sheet view.
$ws.Range("D3").Select()
